$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LAMR")

# Insert two new columns before column D (shifts existing D:K -> F:M).
$ws.Columns("D:E").Insert()

# The inserted columns default to the formatting of the column to their
# left (C); copy number formatting from the (now-shifted) F/G columns so
# the new D/E columns match the data columns instead.
$ws.Range("F7:F102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$ws.Range("G7:G102").Copy()
$ws.Range("E7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the two new quarter columns (D = most recent quarter, E = prior
# quarter) with the newly reported figures.
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 427900
$ws.Range("E8").Value = 418500
$ws.Range("D9").Value = 142100
$ws.Range("E9").Value = 140700
$ws.Range("D10").Value = 285800
$ws.Range("E10").Value = 277800
$ws.Range("D12:E12").Value = "NA"
$ws.Range("D13:E13").Value = 0
$ws.Range("D14:E14").Value = 0
$ws.Range("D15").Value = 58000
$ws.Range("E15").Value = 55100
$ws.Range("D17").Value = 297300
$ws.Range("E17").Value = 290100
$ws.Range("D18").Value = 130600
$ws.Range("E18").Value = 128400
$ws.Range("D20:E20").Value = 200
$ws.Range("D21").Value = 188800
$ws.Range("E21").Value = 183600
$ws.Range("D22").Value = 32400
$ws.Range("E22").Value = 31900
$ws.Range("D23").Value = 98400
$ws.Range("E23").Value = 96700
$ws.Range("D24").Value = 2800
$ws.Range("E24").Value = 2600
$ws.Range("D25:E25").Value = 0
$ws.Range("D26").Value = 95600
$ws.Range("E26").Value = 94100
$ws.Range("D27").Value = 95500
$ws.Range("E27").Value = 94000
$ws.Range("D28:E28").Value = 0
$ws.Range("D29").Value = 100
$ws.Range("E29").Value = "NA"
$ws.Range("D30:E30").Value = 0
$ws.Range("D31:E31").Value = 0
$ws.Range("D32:E32").Value = -200
$ws.Range("D33").Value = 95600
$ws.Range("E33").Value = 94000
$ws.Range("D34:E34").Value = 0
$ws.Range("D35").Value = 95600
$ws.Range("E35").Value = 94000
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 21500
$ws.Range("E41").Value = 10600
$ws.Range("D42:E42").Value = 0
$ws.Range("D43").Value = 235600
$ws.Range("E43").Value = 228000
$ws.Range("D44:E44").Value = 0
$ws.Range("D45").Value = 106600
$ws.Range("E45").Value = 127300
$ws.Range("D46").Value = 363700
$ws.Range("E46").Value = 365900
$ws.Range("D47:E47").Value = 0
$ws.Range("D48").Value = 1295000
$ws.Range("E48").Value = 1193300
$ws.Range("D49").Value = 2834800
$ws.Range("E49").Value = 2511100
$ws.Range("D50:E50").Value = 0
$ws.Range("D51:E51").Value = 0
$ws.Range("D52").Value = 51100
$ws.Range("E52").Value = 54300
$ws.Range("D53:E53").Value = 0
$ws.Range("D54").Value = 4544600
$ws.Range("E54").Value = 4124600
$ws.Range("D57").Value = 21200
$ws.Range("E57").Value = 19500
$ws.Range("D58").Value = 204100
$ws.Range("E58").Value = 26700
$ws.Range("D59").Value = 229700
$ws.Range("E59").Value = 209700
$ws.Range("D60").Value = 455000
$ws.Range("E60").Value = 256000
$ws.Range("D61").Value = 2684600
$ws.Range("E61").Value = 2519200
$ws.Range("D62").Value = 273300
$ws.Range("E62").Value = 252500
$ws.Range("D63:E63").Value = 0
$ws.Range("D64:E64").Value = 0
$ws.Range("D65:E65").Value = 0
$ws.Range("D66").Value = 3412900
$ws.Range("E66").Value = 3027700
$ws.Range("D68:E68").Value = 0
$ws.Range("D69:E69").Value = 0
$ws.Range("D70:E70").Value = 0
$ws.Range("D71:E71").Value = 0
$ws.Range("D72").Value = -695300
$ws.Range("E72").Value = -699300
$ws.Range("D73:E73").Value = 0
$ws.Range("D74:E74").Value = 0
$ws.Range("D75:E75").Value = 0
$ws.Range("D76").Value = 1131800
$ws.Range("E76").Value = 1097000
$ws.Range("D77:E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 95600
$ws.Range("E81").Value = 94000
$ws.Range("D83").Value = 58000
$ws.Range("E83").Value = 55100
$ws.Range("D84:E84").Value = 0
$ws.Range("D85:E85").Value = 0
$ws.Range("D86:E86").Value = 0
$ws.Range("D87:E87").Value = 0
$ws.Range("D88:E88").Value = 0
$ws.Range("D89").Value = 194800
$ws.Range("E89").Value = 154300
$ws.Range("D91").Value = -35500
$ws.Range("E91").Value = -29700
$ws.Range("D92:E92").Value = 0
$ws.Range("D93:E93").Value = 0
$ws.Range("D94").Value = -463800
$ws.Range("E94").Value = -58900
$ws.Range("D96").Value = -91700
$ws.Range("E96").Value = -90300
$ws.Range("D97:E97").Value = 0
$ws.Range("D98:E98").Value = 0
$ws.Range("D99:E99").Value = 0
$ws.Range("D100").Value = 280400
$ws.Range("E100").Value = -104400
$ws.Range("D101").Value = -400
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = 10900
$ws.Range("E102").Value = -9000

Write-Output "Done"
